# Refined metadata to be additional tab
#
# 1) Refresh the "time_taken" timestamps on the existing "data" sheet
#    (re-run of the panel query captured newer per-row timestamps).
# 2) Add a new "metadata" worksheet (after "data") summarising the panel
#    query itself: data_name / data_id / data_version / data_version_created
#    / panel_query_time / panel_get_request.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "data"

# --- 1) refresh time_taken column (F2:F45) on the "data" sheet ---------
$ws1.Range("F2").Value = "2021-10-05 14:21:41.133294"
$ws1.Range("F3").Value = "2021-10-05 14:21:41.133302"
$ws1.Range("F4").Value = "2021-10-05 14:21:41.133306"
$ws1.Range("F5").Value = "2021-10-05 14:21:41.133308"
$ws1.Range("F6").Value = "2021-10-05 14:21:41.133311"
$ws1.Range("F7").Value = "2021-10-05 14:21:41.133314"
$ws1.Range("F8").Value = "2021-10-05 14:21:41.133317"
$ws1.Range("F9").Value = "2021-10-05 14:21:41.133319"
$ws1.Range("F10").Value = "2021-10-05 14:21:41.133323"
$ws1.Range("F11").Value = "2021-10-05 14:21:41.133325"
$ws1.Range("F12").Value = "2021-10-05 14:21:41.133328"
$ws1.Range("F13").Value = "2021-10-05 14:21:41.133331"
$ws1.Range("F14").Value = "2021-10-05 14:21:41.133334"
$ws1.Range("F15").Value = "2021-10-05 14:21:41.133336"
$ws1.Range("F16").Value = "2021-10-05 14:21:41.133339"
$ws1.Range("F17").Value = "2021-10-05 14:21:41.133342"
$ws1.Range("F18").Value = "2021-10-05 14:21:41.133345"
$ws1.Range("F19").Value = "2021-10-05 14:21:41.133347"
$ws1.Range("F20").Value = "2021-10-05 14:21:41.133350"
$ws1.Range("F21").Value = "2021-10-05 14:21:41.133353"
$ws1.Range("F22").Value = "2021-10-05 14:21:41.133356"
$ws1.Range("F23").Value = "2021-10-05 14:21:41.133359"
$ws1.Range("F24").Value = "2021-10-05 14:21:41.133362"
$ws1.Range("F25").Value = "2021-10-05 14:21:41.133365"
$ws1.Range("F26").Value = "2021-10-05 14:21:41.133368"
$ws1.Range("F27").Value = "2021-10-05 14:21:41.133371"
$ws1.Range("F28").Value = "2021-10-05 14:21:41.133373"
$ws1.Range("F29").Value = "2021-10-05 14:21:41.133376"
$ws1.Range("F30").Value = "2021-10-05 14:21:41.133378"
$ws1.Range("F31").Value = "2021-10-05 14:21:41.133381"
$ws1.Range("F32").Value = "2021-10-05 14:21:41.133384"
$ws1.Range("F33").Value = "2021-10-05 14:21:41.133387"
$ws1.Range("F34").Value = "2021-10-05 14:21:41.133390"
$ws1.Range("F35").Value = "2021-10-05 14:21:41.133392"
$ws1.Range("F36").Value = "2021-10-05 14:21:41.133395"
$ws1.Range("F37").Value = "2021-10-05 14:21:41.133398"
$ws1.Range("F38").Value = "2021-10-05 14:21:41.133400"
$ws1.Range("F39").Value = "2021-10-05 14:21:41.133403"
$ws1.Range("F40").Value = "2021-10-05 14:21:41.133406"
$ws1.Range("F41").Value = "2021-10-05 14:21:41.133409"
$ws1.Range("F42").Value = "2021-10-05 14:21:41.133412"
$ws1.Range("F43").Value = "2021-10-05 14:21:41.133415"
$ws1.Range("F44").Value = "2021-10-05 14:21:41.133418"
$ws1.Range("F45").Value = "2021-10-05 14:21:41.133420"

# --- 2) add the "metadata" sheet, positioned right after "data" --------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Header row (bold, bordered, centered -- same look as the "data" header)
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
$col = 2
foreach ($h in $headers) {
    $cell = $ws2.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $col = $col + 1
}

# Data row
$a2 = $ws2.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$ws2.Range("B2").Value = "Multiple monogenic benign skin tumours"
$ws2.Range("C2").Value = 558
$ws2.Range("D2").Value = "'1.12"
$ws2.Range("E2").Value = "2021-03-24T13:44:47.848879Z"
$ws2.Range("F2").Value = "2021-10-05 14:21:41.129841"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/558/?format=json"

# Keep "data" as the active sheet/tab (matches the original activeTab=0)
$ws1.Activate()
$ws1.Range("A1").Select()
